$d = $word.ActiveDocument
$d.Content.Find.Execute("2011-11-25", $true, $false, $false, $false, $false, $true, 1, $false, "2011-11-26", 2)
